# Resource_And_Tracker.xlsx - "Add files via upload" update
#
# Net effect (reconstructed from the OOXML diff):
#  - Tracking sheet: the whole day/session schedule was regenerated starting
#    four days later (Mon 2025-01-06 instead of Thu 2025-01-02), and the
#    weekly pattern changed from "6 sessions/week (Mon-Sat, Sun off)" to
#    "5 sessions/week (Mon-Fri, Sat+Sun off)". Because fewer days/week carry
#    a session, the same 41 sessions (Session-6 .. Session-86, stored as
#    shared-string pairs D/E incrementing by one each workday) now spill
#    into 4 extra rows (through row 56 instead of row 52). The single
#    leftover/odd session at the end (Session-86) is written with no
#    matching E value, same as the original table ended with a lone D value.
#  - Selection/active-cell bookkeeping changed on Check_Points and Tracking.
#
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Tracking sheet: rebuild the date / weekday / session columns
# ---------------------------------------------------------------------
$tracking = $wb.Worksheets.Item("Tracking")

$weekdayNames = @("MON", "TUE", "WED", "THU", "FRI", "SAT", "SUN")

# New schedule starts Monday 2025-01-06 (serial 45663) instead of the old
# Thursday 2025-01-02 (serial 45659).
$startSerial = 45663
$numRows = 55          # rows 2..56
$skipWeekdays = @("SAT", "SUN")

# Sessions are the shared strings "Session-6" (id 124) .. "Session-86" (id 204);
# here we just use the literal text since Excel/the engine resolves the
# shared-string table for us.
$sessionStart = 6
$sessionEnd = 86
$session = $sessionStart

for ($i = 0; $i -lt $numRows; $i++) {
    $row = 2 + $i
    $serial = $startSerial + $i
    $wd = $weekdayNames[$i % 7]

    $tracking.Range("A$row").Value = $row - 1
    $tracking.Range("B$row").Value = $serial
    $tracking.Range("B$row").NumberFormat = "d-mmm"
    $tracking.Range("C$row").Value = $wd

    $remainingSessions = [Math]::Floor(($sessionEnd - $session) / 1) + 1
    $isSkipDay = $skipWeekdays -contains $wd

    $placeHere = $false
    if ($session -le $sessionEnd) {
        if (-not $isSkipDay) {
            $placeHere = $true
        } elseif ($remainingSessions -eq 1) {
            # lone leftover session spills onto the next row even if it
            # would normally be a day off
            $placeHere = $true
        }
    }

    if ($placeHere) {
        $tracking.Range("D$row").Value = "Session-$session"
        $session = $session + 1
        if ($session -le $sessionEnd) {
            $tracking.Range("E$row").Value = "Session-$session"
            $session = $session + 1
        }
    }
}

# Sheet view bookkeeping for Tracking
try {
    $tracking.Application.ActiveWindow.ScrollRow = 31
} catch {
}
$tracking.Range("F53").Select()

# ---------------------------------------------------------------------
# 2) Check_Points sheet: selection moved to F26
# ---------------------------------------------------------------------
$checkPoints = $wb.Worksheets.Item("Check_Points")
$checkPoints.Range("F26").Select()

# ---------------------------------------------------------------------
# 3) Topics sheet: no data/selection changes, leave as-is
# ---------------------------------------------------------------------

# Re-activate Tracking (it was the active sheet/tab before and after)
$tracking.Activate()
$tracking.Range("F53").Select()
